$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46060 -> 46061) for every data row (rows 2 through 438).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 438) {
    $lastRow = 438
}

$ws.Range("C2:C$lastRow").Value = 46061
